# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the existing header cells (A1:AC1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-40 - same record for every player row
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 101   # column AD
    $ws.Cells.Item($r, 31).Value = 59    # column AE
    $ws.Cells.Item($r, 32).Value = 1     # column AF
}
